# Corrections following third round of review
#
# The "Materials" sheet had an erroneous "subgenus" column (column AS,
# header "subgenus" / data "${subgenus}") that duplicated/was not wanted
# alongside "genus". Remove that entire column, which shifts every
# subsequent column one position to the left and lets the now-unused
# "subgenus" / "${subgenus}" shared strings drop out of the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets("Materials")
$ws.Columns("AS:AS").Delete()
